# Capitalize the first letter of several PHENOTYPE values in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of old text -> new text (capitalize first letter only)
$replacements = @{
    "endometriosis"                 = "Endometriosis"
    "polycystic ovary syndrome"     = "Polycystic ovary syndrome"
    "recurrent spontaneous abortion" = "Recurrent spontaneous abortion"
    "oligoasthenoteratozoospermia"  = "Oligoasthenoteratozoospermia"
    "non-obstructive azoospermia"   = "Non-obstructive azoospermia"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
